$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple numeric updates (same country, same row) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value2 = 1137494
$ws.Range("C4").Value2 = 6464
$ws.Range("D4").Value2 = 162100
$ws.Range("E4").Value2 = 909123
$ws.Range("G4").Value2 = 518
$ws.Range("H4").Value2 = 66271

# Row 6: Italia
$ws.Range("B6").Value2 = 209328
$ws.Range("C6").Value2 = 1900
$ws.Range("D6").Value2 = 79914
$ws.Range("E6").Value2 = 100704
$ws.Range("F6").Value2 = 1539
$ws.Range("G6").Value2 = 474
$ws.Range("H6").Value2 = 28710

# Row 7: Reino Unido
$ws.Range("B7").Value2 = 182260
$ws.Range("C7").Value2 = 4806
$ws.Range("E7").Value2 = 153785
$ws.Range("G7").Value2 = 621
$ws.Range("H7").Value2 = 28131

# Row 37: Rumania
$ws.Range("E37").Value2 = 7414
$ws.Range("G37").Value2 = 27
$ws.Range("H37").Value2 = 771

# Row 60: Luxemburgo
$ws.Range("B60").Value2 = 3812
$ws.Range("C60").Value2 = 10
$ws.Range("D60").Value2 = 3318
$ws.Range("E60").Value2 = 402

# --- Re-rank of Irak / Nigeria / Ghana (rows 69-71) ---
# Irak's case count jumped ahead of Nigeria and Ghana, so it is now listed
# right after Armenia, pushing Nigeria and Ghana each one row down.

# Row 69 now holds Irak's (new) figures
$ws.Range("A69").Value2 = "Irak"
$ws.Range("B69").Value2 = 2219
$ws.Range("C69").Value2 = 66
$ws.Range("D69").Value2 = 1473
$ws.Range("E69").Value2 = 651
$ws.Range("F69").Value2 = 0
$ws.Range("G69").Value2 = 1
$ws.Range("H69").Value2 = 95

# Row 70 now holds Nigeria's (unchanged) figures
$ws.Range("A70").Value2 = "Nigeria"
$ws.Range("B70").Value2 = 2170
$ws.Range("C70").Value2 = 0
$ws.Range("D70").Value2 = 351
$ws.Range("E70").Value2 = 1751
$ws.Range("F70").Value2 = 2
$ws.Range("G70").Value2 = 0
$ws.Range("H70").Value2 = 68

# Row 71 now holds Ghana's (unchanged) figures
$ws.Range("A71").Value2 = "Ghana"
$ws.Range("B71").Value2 = 2169
$ws.Range("C71").Value2 = 95
$ws.Range("D71").Value2 = 229
$ws.Range("E71").Value2 = 1922
$ws.Range("F71").Value2 = 4
$ws.Range("G71").Value2 = 1
$ws.Range("H71").Value2 = 18
